$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (rows 401-403), shifting
# the existing rows 401-419 down to 404-422, matching the weekly update.
$ws.Rows("401:403").Insert()

# Row 401
$ws.Range("A401").Value = 11
$ws.Range("B401").Value = "Vega Monumental Concepción"
$ws.Range("C401").Value = "Bíobío"
$ws.Range("D401").Value = 44516
$ws.Range("E401").Value = 8
$ws.Range("F401").Value = 100112033
$ws.Range("G401").Value = "Lechuga"
$ws.Range("H401").Value = "Conconina(o)"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 450
$ws.Range("K401").Value = 4000
$ws.Range("L401").Value = 4500
$ws.Range("M401").Value = 4222
$ws.Range("N401").Value = "`$/caja 10 unidades"
$ws.Range("O401").Value = "Región Metropolitana"
$ws.Range("P401").Value = 422
$ws.Range("Q401").Value = 10
$ws.Range("R401").Value = "Hortaliza"

# Row 402
$ws.Range("A402").Value = 11
$ws.Range("B402").Value = "Vega Monumental Concepción"
$ws.Range("C402").Value = "Bíobío"
$ws.Range("D402").Value = 44516
$ws.Range("E402").Value = 8
$ws.Range("F402").Value = 100112033
$ws.Range("G402").Value = "Lechuga"
$ws.Range("H402").Value = "Escarola"
$ws.Range("I402").Value = "Primera"
$ws.Range("J402").Value = 450
$ws.Range("K402").Value = 5000
$ws.Range("L402").Value = 5500
$ws.Range("M402").Value = 5222
$ws.Range("N402").Value = "`$/caja 15 unidades"
$ws.Range("O402").Value = "Región Metropolitana"
$ws.Range("P402").Value = 348
$ws.Range("Q402").Value = 15
$ws.Range("R402").Value = "Hortaliza"

# Row 403
$ws.Range("A403").Value = 11
$ws.Range("B403").Value = "Vega Monumental Concepción"
$ws.Range("C403").Value = "Bíobío"
$ws.Range("D403").Value = 44516
$ws.Range("E403").Value = 8
$ws.Range("F403").Value = 100112033
$ws.Range("G403").Value = "Lechuga"
$ws.Range("H403").Value = "Marina"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 250
$ws.Range("K403").Value = 4500
$ws.Range("L403").Value = 5000
$ws.Range("M403").Value = 4700
$ws.Range("N403").Value = "`$/caja 15 unidades"
$ws.Range("O403").Value = "Región Metropolitana"
$ws.Range("P403").Value = 313
$ws.Range("Q403").Value = 15
$ws.Range("R403").Value = "Hortaliza"
